# Update the BNF grammar rules for operator precedence (<оператор_1> .. <оператор_9>)
# and the cached page-number field result in the footer.

$d = $word.ActiveDocument

# Belt-and-braces: make sure Word doesn't "smarten" straight quotes into
# curly quotes while we touch the grammar text below (some of the BNF
# rules contain a literal "||" token that must stay as straight quotes).
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

function Replace-Exact {
    param(
        [string]$old,
        [string]$new
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $old"
    }
    # Setting .Text directly (instead of passing replacement text to
    # Find.Execute) avoids Word's smart-quote autoformat substitution,
    # preserving straight double quotes in the grammar text.
    $rng.Text = $new
}

Replace-Exact " ::= (<оператор_1>|<оператор>{,<оператор_1>})" " ::= <оператор_1>{,<оператор_1>}"
Replace-Exact " ::= (<оператор_2>|<объект>(=|+=|-=|*=|/=|%=)<оператор_2>)" " ::= [<объект>(=|+=|-=|*=|/=|%=)]<оператор_2>"
Replace-Exact " ::= (<оператор_3>|<оператор_2>`"||`"<оператор_3>)" " ::= <оператор_3>[`"||`"<оператор_3>]"
Replace-Exact " ::= (<оператор_4>|<оператор_3>&&<оператор_4>)" " ::= <оператор_4>[&&<оператор_4>]"
Replace-Exact " ::= (<оператор_5>|<оператор_4>(==|!=)<оператор_5>)" " ::= <оператор_5>[(==|!=)<оператор_5>]"
Replace-Exact " ::= (<оператор_6>|<оператор_5>(<|<=|>|>=)<оператор_6>)" " ::= <оператор_6>[(<|<=|>|>=)<оператор_6>]"
Replace-Exact " ::= (<оператор_7>|<оператор_6>(+|-)<оператор_7>)" " ::= <оператор_7>[(+|-)<оператор_7>]"
Replace-Exact " ::= (<оператор_8>|<оператор_7>(*|/|%)<оператор_8>)" " ::= <оператор_8>[(*|/|%)<оператор_8>]"
Replace-Exact " ::= (<оператор_9>|(++|--)<объект>|(+|-|!|*)<оператор_9>|&<объект>)" " ::= (<оператор_9>[(+|-|!|*)<оператор_9>]|(++|--|&)<объект>)"

# Footer page-number cached field result: 4 -> 3
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item(1)
    if ($ftr.Exists) {
        foreach ($fld in $ftr.Range.Fields) {
            if ($fld.Result.Text.Trim() -eq "4") {
                $fld.Result.Text = "3"
            }
        }
    }
}
